$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Hunk 1 -------------------------------------------------------------
# Collapse the 4 paragraphs:
#   ".Uitbreiding van ondergrondse tunnel Vault -> Graveyard {New end}"
#   "Monsters {Special item om het canon te gebruiken}, One way tunnel"
#   ".{Special item om canon op schip te gebruiken} "
#   ".{Plek in de staat die er dan kapot uitziet door schip cannon}"
# into a single new paragraph about the "Blown up housse" secret tunnel.
$p1Start = $d.Paragraphs(3).Range.Start
$p1End = $d.Paragraphs(6).Range.End
$r1 = $d.Range($p1Start, $p1End)
$xml1 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>.{</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Blown up </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>housse</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> -&gt; End 5 </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>serect</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> tunnel</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}</w:t></w:r>' +
  '</w:p>'
$r1.InsertXML($xml1)

# --- Hunk 2 -------------------------------------------------------------
# ".Quanity box for the shops" -> ".Let the Carrier ride to a abonded
# {Marine base -> sub or something like abdonded}"
# (hunk 1 collapsed 4 paragraphs into 1, so indices shift down by 3)
# The paragraph mark / w:p attributes are untouched by this hunk, so keep
# the original w:rsidR="00F354DD" w:rsidRDefault="00F354DD" around.
$p2 = $d.Paragraphs(5).Range
$xml2 = '<w:p ' + $wns + ' w:rsidR="00F354DD" w:rsidRDefault="00F354DD">' +
  '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>.Let</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> the Carrier ride to a </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>abonded</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> {Marine base </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> sub or something like </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>abdonded</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}</w:t></w:r>' +
  '</w:p>'
$p2.InsertXML($xml2)

# --- Hunk 3 -------------------------------------------------------------
# ".Let the Carrier ride to a abonded {Marine base -> sub or something
# like abdonded}" -> ".NPC Interact systeem" (bookmark _GoBack retained)
# Keep the original w:p attributes (untouched by this hunk).
$p3 = $d.Paragraphs(6).Range
$xml3 = '<w:p ' + $wns + ' w:rsidR="009F6113" w:rsidRPr="00F354DD" w:rsidRDefault="009F6113">' +
  '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">.NPC </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Interact systeem</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$p3.InsertXML($xml3)
